$wb = $excel.ActiveWorkbook

# --- AMSIN sheet: append rows 47-50 ---
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$amsinRows = @(
    @("2022-09-12", 44816.62676319444, "ecstest167", 105, 99, 6, 4.87),
    @("2022-09-16", 44820.63874326389, "fstcyc167", 105, 105, 0, 3.26),
    @("2022-09-19", 44823.66180796296, "scndcycle167", 105, 94, 11, 5.15),
    @("2022-09-20", 44824.3776330787, "finalrun167", 105, 100, 5, 4.01)
)

$r = 47
foreach ($row in $amsinRows) {
    $wsAmsin.Cells.Item($r, 1).Value = $row[0]
    $wsAmsin.Cells.Item($r, 2).Value = $row[1]
    $wsAmsin.Cells.Item($r, 3).Value = $row[2]
    $wsAmsin.Cells.Item($r, 4).Value = $row[3]
    $wsAmsin.Cells.Item($r, 5).Value = $row[4]
    $wsAmsin.Cells.Item($r, 6).Value = $row[5]
    $wsAmsin.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# --- BETA sheet: append row 24 ---
$wsBeta = $wb.Worksheets.Item("BETA")
$wsBeta.Cells.Item(24, 1).Value = "2022-09-20"
$wsBeta.Cells.Item(24, 2).Value = 44824.52734270386
$wsBeta.Cells.Item(24, 3).Value = "beta167"
$wsBeta.Cells.Item(24, 4).Value = 105
$wsBeta.Cells.Item(24, 5).Value = 105
$wsBeta.Cells.Item(24, 6).Value = 0
$wsBeta.Cells.Item(24, 7).Value = 2.75

# --- AMS sheet: fix row 22 time value ---
$wsAms = $wb.Worksheets.Item("AMS")
$wsAms.Cells.Item(22, 2).Value = 44797.92717719908
